# Update Betfair Back/Lay odds in Jogos_do_Dia_Betfair_Back_Lay_2025-11-19.xlsx
# (daily refresh of odds values for rows 2-11, columns F:AO)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.74
$ws.Range("H2").Value = 2.76
$ws.Range("I2").Value = 3.15
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 2.68
$ws.Range("O2").Value = 1.47
$ws.Range("P2").Value = 1.57
$ws.Range("Q2").Value = 2.42
$ws.Range("R2").Value = 1.21
$ws.Range("S2").Value = 5
$ws.Range("U2").Value = 1.84
$ws.Range("X2").Value = 9.800000000000001
$ws.Range("Y2").Value = 9.4
$ws.Range("Z2").Value = 19
$ws.Range("AB2").Value = 9.4
$ws.Range("AC2").Value = 7.4
$ws.Range("AD2").Value = 14
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 19
$ws.Range("AH2").Value = 65
$ws.Range("AJ2").Value = 200
$ws.Range("AK2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3
$ws.Range("F3").Value = 2.56
$ws.Range("G3").Value = 2.86
$ws.Range("H3").Value = 2.94
$ws.Range("I3").Value = 3.4
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 3.4
$ws.Range("N3").Value = 2.92
$ws.Range("O3").Value = 1.43
$ws.Range("P3").Value = 1.54
$ws.Range("S3").Value = 3.15
$ws.Range("V3").Value = 1.42
$ws.Range("W3").Value = 1.53
$ws.Range("X3").Value = 11
$ws.Range("Y3").Value = 10.5
$ws.Range("Z3").Value = 22
$ws.Range("AA3").Value = 900
$ws.Range("AB3").Value = 9.4
$ws.Range("AC3").Value = 7.6
$ws.Range("AD3").Value = 14.5
$ws.Range("AF3").Value = 17.5
$ws.Range("AG3").Value = 13
$ws.Range("AH3").Value = 21
$ws.Range("AJ3").Value = 200
$ws.Range("AK3").Value = 200
$ws.Range("AL3").Value = 260
$ws.Range("AN3").Value = 42
$ws.Range("AO3").Value = 55

# Row 4
$ws.Range("F4").Value = 1.31
$ws.Range("G4").Value = 1.33
$ws.Range("H4").Value = 13
$ws.Range("I4").Value = 14
$ws.Range("J4").Value = 5.9
$ws.Range("K4").Value = 6.2
$ws.Range("N4").Value = 4.3
$ws.Range("O4").Value = 1.27
$ws.Range("P4").Value = 2.16
$ws.Range("Q4").Value = 1.81
$ws.Range("R4").Value = 1.43
$ws.Range("S4").Value = 3.1
$ws.Range("U4").Value = 1.7
$ws.Range("V4").Value = 1.07
$ws.Range("W4").Value = 4
$ws.Range("X4").Value = 20
$ws.Range("Y4").Value = 38
$ws.Range("Z4").Value = 130
$ws.Range("AA4").Value = 700
$ws.Range("AC4").Value = 13.5
$ws.Range("AD4").Value = 48
$ws.Range("AE4").Value = 270
$ws.Range("AF4").Value = 7
$ws.Range("AH4").Value = 38
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 9.6
$ws.Range("AK4").Value = 15
$ws.Range("AL4").Value = 46
$ws.Range("AN4").Value = 5.9
$ws.Range("AO4").Value = 420

# Row 5
$ws.Range("G5").Value = 1.99
$ws.Range("H5").Value = 4.9
$ws.Range("I5").Value = 5.6
$ws.Range("Q5").Value = 2.36

# Row 6
$ws.Range("F6").Value = 2.4
$ws.Range("G6").Value = 2.92
$ws.Range("H6").Value = 2.72
$ws.Range("I6").Value = 3.45
$ws.Range("J6").Value = 3.1
$ws.Range("K6").Value = 3.95
$ws.Range("L6").Value = 1.34
$ws.Range("P6").Value = 1.89
$ws.Range("Y6").Value = 1000
$ws.Range("AA6").Value = 1000

# Row 7
$ws.Range("F7").Value = 1.99
$ws.Range("I7").Value = 5.2
$ws.Range("J7").Value = 3.15
$ws.Range("L7").Value = 1.55
$ws.Range("N7").Value = 2.66
$ws.Range("O7").Value = 1.51
$ws.Range("P7").Value = 1.56
$ws.Range("Q7").Value = 2.5
$ws.Range("S7").Value = 5.3
$ws.Range("T7").Value = 2.14
$ws.Range("V7").Value = 1.25
$ws.Range("AA7").Value = 160
$ws.Range("AB7").Value = 15
$ws.Range("AC7").Value = 17
$ws.Range("AD7").Value = 25
$ws.Range("AF7").Value = 42
$ws.Range("AI7").Value = 470

# Row 8
$ws.Range("I8").Value = 2.16
$ws.Range("J8").Value = 3.3
$ws.Range("K8").Value = 3.4
$ws.Range("N8").Value = 2.8
$ws.Range("O8").Value = 1.52
$ws.Range("Q8").Value = 2.56
$ws.Range("S8").Value = 5.3
$ws.Range("U8").Value = 1.79
$ws.Range("V8").Value = 1.86
$ws.Range("Y8").Value = 7
$ws.Range("AM8").Value = 190

# Row 9
$ws.Range("F9").Value = 2.24
$ws.Range("G9").Value = 2.3
$ws.Range("H9").Value = 3.55
$ws.Range("I9").Value = 3.7
$ws.Range("J9").Value = 3.5
$ws.Range("K9").Value = 3.6
$ws.Range("N9").Value = 3.55
$ws.Range("P9").Value = 1.86
$ws.Range("R9").Value = 1.33
$ws.Range("T9").Value = 1.85
$ws.Range("U9").Value = 2.06
$ws.Range("W9").Value = 1.77
$ws.Range("X9").Value = 13
$ws.Range("Z9").Value = 25
$ws.Range("AA9").Value = 70
$ws.Range("AC9").Value = 7.6
$ws.Range("AD9").Value = 15
$ws.Range("AE9").Value = 46
$ws.Range("AF9").Value = 13.5
$ws.Range("AG9").Value = 11
$ws.Range("AH9").Value = 18.5
$ws.Range("AI9").Value = 60
$ws.Range("AJ9").Value = 29
$ws.Range("AK9").Value = 25
$ws.Range("AL9").Value = 42
$ws.Range("AM9").Value = 110
$ws.Range("AN9").Value = 20
$ws.Range("AO9").Value = 48

# Row 10
$ws.Range("F10").Value = 2.3
$ws.Range("H10").Value = 3.45
$ws.Range("K10").Value = 3.55
$ws.Range("N10").Value = 3.4
$ws.Range("P10").Value = 1.84
$ws.Range("R10").Value = 1.32
$ws.Range("T10").Value = 1.88
$ws.Range("U10").Value = 2.06
$ws.Range("W10").Value = 1.74
$ws.Range("Y10").Value = 13
$ws.Range("Z10").Value = 24
$ws.Range("AA10").Value = 65
$ws.Range("AB10").Value = 9.6
$ws.Range("AD10").Value = 15.5
$ws.Range("AE10").Value = 44
$ws.Range("AH10").Value = 20
$ws.Range("AJ10").Value = 32
$ws.Range("AL10").Value = 42
$ws.Range("AM10").Value = 120
$ws.Range("AN10").Value = 23
$ws.Range("AO10").Value = 50

# Row 11
$ws.Range("F11").Value = 2.3
$ws.Range("G11").Value = 2.78
$ws.Range("H11").Value = 2.66
$ws.Range("I11").Value = 3.25
$ws.Range("J11").Value = 3.25
$ws.Range("K11").Value = 4.5
$ws.Range("P11").Value = 1.84
$ws.Range("Q11").Value = 1.78
$ws.Range("S11").Value = 3.15
$ws.Range("V11").Value = 1.44
$ws.Range("W11").Value = 1.56
$ws.Range("AB11").Value = 1000
$ws.Range("AC11").Value = 11.5
$ws.Range("AG11").Value = 1000
$ws.Range("AM11").Value = 1000
